$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.11"
$ws.Range("E2").Value = "'6.34%"
$ws.Range("D3").Value = "'31.99"
$ws.Range("E3").Value = "'9.13%"
$ws.Range("D4").Value = "'5.295"
$ws.Range("E4").Value = "'4.26%"
$ws.Range("D5").Value = "'0.07440"
$ws.Range("E5").Value = "'10.68%"
$ws.Range("D6").Value = "'7.838"
$ws.Range("E6").Value = "'6.99%"
$ws.Range("D7").Value = "'3.774"
$ws.Range("E7").Value = "'9.74%"
$ws.Range("D8").Value = "'1.477"
$ws.Range("E8").Value = "'6.56%"
$ws.Range("D9").Value = "'0.9144"
$ws.Range("E9").Value = "'1.46%"
$ws.Range("D10").Value = "'0.01761"
$ws.Range("E10").Value = "'2,620.00%"
$ws.Range("E11").Value = "'7.36%"
$ws.Range("D12").Value = "'0.07662"
$ws.Range("E12").Value = "'11.31%"
$ws.Range("D13").Value = "'0.08154"
$ws.Range("E13").Value = "'6.83%"
$ws.Range("D14").Value = "'0.03040"
$ws.Range("E14").Value = "'4.01%"
$ws.Range("D15").Value = "'0.09943"
$ws.Range("E15").Value = "'10.60%"
$ws.Range("D16").Value = "'0.001510"
$ws.Range("E16").Value = "'-5.16%"
$ws.Range("D17").Value = "'0.04572"
$ws.Range("E17").Value = "'1.90%"
$ws.Range("D18").Value = "'0.006259"
$ws.Range("E18").Value = "'-4.30%"
$ws.Range("D19").Value = "'3.485"
$ws.Range("E19").Value = "'1.03%"
$ws.Range("D20").Value = "'2.229"
$ws.Range("E20").Value = "'-0.11%"
$ws.Range("D21").Value = "'0.3301"
$ws.Range("E21").Value = "'2.99%"
$ws.Range("D22").Value = "'0.1348"
$ws.Range("E22").Value = "'2.24%"
$ws.Range("D23").Value = "'4.495"
$ws.Range("E23").Value = "'11.32%"
$ws.Range("D24").Value = "'0.1644"
$ws.Range("E24").Value = "'4.09%"
$ws.Range("D25").Value = "'0.001218"
$ws.Range("E25").Value = "'1.41%"
$ws.Range("D26").Value = "'0.004398"
$ws.Range("E26").Value = "'0.60%"
$ws.Range("D27").Value = "'0.0001402"
$ws.Range("E27").Value = "'20.02%"
$ws.Range("D28").Value = "'0.0001748"
$ws.Range("E28").Value = "'8.07%"
$ws.Range("D40").Value = "'0.04518"
$ws.Range("E40").Value = "'6.52%"
$ws.Range("D41").Value = "'0.007069"
$ws.Range("E41").Value = "'4.15%"
$ws.Range("D42").Value = "'0.1344"
$ws.Range("E42").Value = "'8.45%"
$ws.Range("D43").Value = "'0.002243"
$ws.Range("E43").Value = "'2.12%"
$ws.Range("D44").Value = "'0.01358"
$ws.Range("E44").Value = "'17.24%"
$ws.Range("D45").Value = "'0.00006212"
$ws.Range("E45").Value = "'8.54%"
$ws.Range("D46").Value = "'0.7082"
$ws.Range("E46").Value = "'-63.33%"
$ws.Range("D47").Value = "'0.01305"
$ws.Range("E47").Value = "'-13.02%"
